# Applies the "removed Hasith Dewmina from students -> removed_students,
# reset password sheet" edit described by the target diff.

$wb = $excel.ActiveWorkbook

# Helper sheets
$wsStudents  = $wb.Worksheets.Item("students")
$wsRemoved   = $wb.Worksheets.Item("removed_students")
$wsPswd      = $wb.Worksheets.Item("student_pswd")

# ---------------------------------------------------------------------
# Sheet "students": the sample records (rows 2-4) are removed, the
# header in A1 is renamed from "student_admission" to "Index No. ",
# and the student-count cell (J3) drops to 0.
# ---------------------------------------------------------------------
$wsStudents.Range("A1").Value = "Index No. "
$wsStudents.Range("A2:H4").ClearContents()
$wsStudents.Range("J3").Value = 0
$wsStudents.Columns("A").AutoFit()
$wsStudents.Range("F30").Select()

# ---------------------------------------------------------------------
# Sheet "removed_students": one record (Hasith Dewmina) is recorded as
# removed, and the running total in J3 becomes 1.
# ---------------------------------------------------------------------
$wsRemoved.Range("A3:H3").ClearContents()

$wsRemoved.Range("A2").Value = 1
$wsRemoved.Range("B2").Value = "Hasith"
$wsRemoved.Range("C2").Value = "Dewmina"

$wsRemoved.Range("D2").NumberFormat = "@"
$wsRemoved.Range("D2").Value = "76678687"
$wsRemoved.Range("D2").Style = "Normal"

$wsRemoved.Range("E2").Value = "M"
$wsRemoved.Range("F2").Value = "Lesli Kumara, Gangani Madawala"

$wsRemoved.Range("G2").NumberFormat = "@"
$wsRemoved.Range("G2").Value = "874857498379"
$wsRemoved.Range("G2").Style = "Normal"

$wsRemoved.Range("H2").Value = "Ruwi, Muscat, Oman"
$wsRemoved.Range("J2").Value = "Num Students"
$wsRemoved.Range("J3").Value = 1

$wsRemoved.Range("J3").Select()

# ---------------------------------------------------------------------
# Sheet "student_pswd": the stored credentials are cleared out, the
# "Name" column header becomes "First Name", and the row counter (G6)
# drops to 0.
# ---------------------------------------------------------------------
$wsPswd.Range("C1").Value = "First Name"
$wsPswd.Range("A2:C4").ClearContents()
$wsPswd.Range("G6").Value = 0
$wsPswd.Columns("C").AutoFit()
$wsPswd.Range("G6").Select()
